$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,8).Value = 0.1053200318617804

$ws.Cells.Item(3,2).Value = 0.06502862335842799
$ws.Cells.Item(3,8).Value = 0.1703486552202084

$ws.Cells.Item(4,2).Value = 0.05777295219089382
$ws.Cells.Item(4,8).Value = 0.1630929840526742

$ws.Cells.Item(5,2).Value = 0.04003746439820822
$ws.Cells.Item(5,3).Value = 0.005433154263989929
$ws.Cells.Item(5,4).Value = 6.903545396162682
$ws.Cells.Item(5,5).Value = 0.07162190113110609
$ws.Cells.Item(5,6).Value = 0.02933344314718372
$ws.Cells.Item(5,7).Value = 0.05074148564923312
$ws.Cells.Item(5,8).Value = 0.1453574962599886

$ws.Cells.Item(6,2).Value = 0.02668104829557686
$ws.Cells.Item(6,3).Value = 0.004452646875450709
$ws.Cells.Item(6,4).Value = 2.87173876220049
$ws.Cells.Item(6,5).Value = 0.03367168509957227
$ws.Cells.Item(6,6).Value = 0.01790843688379194
$ws.Cells.Item(6,7).Value = 0.0354536597073619
$ws.Cells.Item(6,8).Value = 0.1320010801573573

$ws.Cells.Item(7,2).Value = 0.01905425749260545
$ws.Cells.Item(7,8).Value = 0.1243742893543858
$ws.Cells.Item(7,3).ClearContents()
$ws.Cells.Item(7,4).ClearContents()
$ws.Cells.Item(7,5).ClearContents()
$ws.Cells.Item(7,6).ClearContents()
$ws.Cells.Item(7,7).ClearContents()

$ws.Cells.Item(8,2).Value = 0.01736560068711328
$ws.Cells.Item(8,3).Value = 0.004163947991993477
$ws.Cells.Item(8,4).Value = 1.466885868502567
$ws.Cells.Item(8,5).Value = 0.0172446340642596
$ws.Cells.Item(8,6).Value = 0.009142926149439599
$ws.Cells.Item(8,7).Value = 0.02558827522478653
$ws.Cells.Item(8,8).Value = 0.1226856325488937

$ws.Cells.Item(9,2).Value = 0.01732275822753888
$ws.Cells.Item(9,3).Value = 0.003113653236048236
$ws.Cells.Item(9,4).Value = 1.610124797313698
$ws.Cells.Item(9,5).Value = 0.0159669773847593
$ws.Cells.Item(9,6).Value = 0.01120881527836379
$ws.Cells.Item(9,7).Value = 0.02343670117671362
$ws.Cells.Item(9,8).Value = 0.1226427900893193

$ws.Cells.Item(10,2).Value = 0.0162799505771167
$ws.Cells.Item(10,3).Value = 0.003428943804904028
$ws.Cells.Item(10,4).Value = 1.562044403543934
$ws.Cells.Item(10,5).Value = 0.01730849734998471
$ws.Cells.Item(10,6).Value = 0.009556342929723413
$ws.Cells.Item(10,7).Value = 0.02300355822450985
$ws.Cells.Item(10,8).Value = 0.1215999824388971

$ws.Cells.Item(11,2).Value = 0.02862222822824241
$ws.Cells.Item(11,8).Value = 0.1339422600900228

$ws.Cells.Item(12,2).Value = 0.04430108711942315
$ws.Cells.Item(12,8).Value = 0.1496211189812035

$ws.Cells.Item(13,2).Value = 0.05520838535795865
$ws.Cells.Item(13,8).Value = 0.160528417219739

$ws.Cells.Item(14,2).Value = 0.06313946483995493
$ws.Cells.Item(14,8).Value = 0.1684594967017353

$ws.Cells.Item(15,2).Value = 0.06534465239163116
$ws.Cells.Item(15,8).Value = 0.1706646842534116

$ws.Cells.Item(16,2).Value = 0.06948924553581119
$ws.Cells.Item(16,8).Value = 0.1748092773975916

$ws.Cells.Item(17,2).Value = 0.07183312015261908
$ws.Cells.Item(17,8).Value = 0.1771531520143995

$ws.Cells.Item(18,2).Value = -0.1053200318617804

$ws.Cells.Item(19,2).Value = 0.07568799403071995
$ws.Cells.Item(19,8).Value = 0.1810080258925003

$ws.Cells.Item(20,2).Value = 0.07696599933625284
$ws.Cells.Item(20,8).Value = 0.1822860311980332

$ws.Cells.Item(21,2).Value = 0.07958835759029037
$ws.Cells.Item(21,8).Value = 0.1849083894520708

$ws.Cells.Item(22,2).Value = 0.08066186028953697
$ws.Cells.Item(22,3).Value = 0.009749602774917389
$ws.Cells.Item(22,4).Value = 14.9528240895545
$ws.Cells.Item(22,5).Value = 0.05442039328382266
$ws.Cells.Item(22,6).Value = 0.06149311110204266
$ws.Cells.Item(22,7).Value = 0.0998306094770315
$ws.Cells.Item(22,8).Value = 0.1859818921513174

$ws.Cells.Item(23,2).Value = 0.08147063306007334
$ws.Cells.Item(23,8).Value = 0.1867906649218537

$ws.Cells.Item(24,2).Value = 0.07791512637333038
$ws.Cells.Item(24,3).Value = 0.009884813500305314
$ws.Cells.Item(24,4).Value = 12.73343963955979
$ws.Cells.Item(24,5).Value = 0.05946162548731106
$ws.Cells.Item(24,6).Value = 0.05848277880592048
$ws.Cells.Item(24,7).Value = 0.09734747394074011
$ws.Cells.Item(24,8).Value = 0.1832351582351108

$ws.Cells.Item(25,2).Value = 0.07538410366498602
$ws.Cells.Item(25,3).Value = 0.01044901290117689
$ws.Cells.Item(25,4).Value = 10.92792185528193
$ws.Cells.Item(25,5).Value = 0.08118253572269221
$ws.Cells.Item(25,6).Value = 0.054797339215146
$ws.Cells.Item(25,7).Value = 0.09597086811482593
$ws.Cells.Item(25,8).Value = 0.1807041355267664

$ws.Cells.Item(26,2).Value = 0.07672912400592764
$ws.Cells.Item(26,3).Value = 0.0095112734571934
$ws.Cells.Item(26,4).Value = -615458881264.1466
$ws.Cells.Item(26,5).Value = 0.06873298698064179
$ws.Cells.Item(26,6).Value = 0.05802678591276553
$ws.Cells.Item(26,7).Value = 0.09543146209908965
$ws.Cells.Item(26,8).Value = 0.182049155867708

$ws.Cells.Item(27,2).Value = 0.07461329212130874
$ws.Cells.Item(27,3).Value = 0.00947194696835187
$ws.Cells.Item(27,4).Value = 10.08078913534554
$ws.Cells.Item(27,5).Value = 0.07725968827376228
$ws.Cells.Item(27,6).Value = 0.05597423592445557
$ws.Cells.Item(27,7).Value = 0.09325234831816163
$ws.Cells.Item(27,8).Value = 0.1799333239830891

$ws.Cells.Item(28,2).Value = 0.0716448291081884
$ws.Cells.Item(28,3).Value = 0.00929386048126217
$ws.Cells.Item(28,4).Value = 9.783011301333788
$ws.Cells.Item(28,5).Value = 0.1032664667766986
$ws.Cells.Item(28,6).Value = 0.05334753270758397
$ws.Cells.Item(28,7).Value = 0.08994212550879278
$ws.Cells.Item(28,8).Value = 0.1769648609699688

$ws.Cells.Item(29,2).Value = 0.01837697764680676
$ws.Cells.Item(29,3).Value = 0.003372328061941301
$ws.Cells.Item(29,4).Value = 1.910192479225142
$ws.Cells.Item(29,5).Value = 0.01332468734440235
$ws.Cells.Item(29,6).Value = 0.01172144556022861
$ws.Cells.Item(29,7).Value = 0.02503250973338464
$ws.Cells.Item(29,8).Value = 0.1236970095085872

